$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medicos")

$ws.Range("E5").Value = "{{Medico.EstadoId}}"
$ws.Range("E7").Value = "{{Medico.CiudadId}}"

$ws.Range("E7").Select()
